$d = $word.ActiveDocument

# --- Paragraph 2 ("Action Yapısı : ...") : append two new runs ---
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertAfter("Actions redux ile ilk iletişimi kurdugumuz noktadır.")
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertAfter(" ")

# --- Paragraph 3 (the bookmark paragraph) : insert the new sentences
#     right before the _GoBack bookmark, inside the SAME paragraph.
#     Build the whole run of text first (so every insertion lands in
#     plain/non-bold formatting), then go back and apply Bold to just
#     the two specific words ("Reducer" and "Store") using position
#     arithmetic computed from the known run lengths - this avoids the
#     "new text inherits the bold of what precedes it" trap that
#     InsertBefore/Font.Bold ordering would otherwise hit.

$run1 = "Reducer ve Store:Action yapısı Reducer ile kontrol altına alınıyor."
$run2 = "Reducer"
$run3 = " ilgili alanın içerip içermediğine bakıyor."
$run4 = "Bu yapının tamamamını içeren sistemede "
$run5 = "Store"
$run6 = " diyoruz."

$fullText = $run1 + $run2 + $run3 + $run4 + $run5 + $run6

$bookmark = $d.Bookmarks.Item("_GoBack")
$paraStart = $bookmark.Range.Start
$d.Range($paraStart, $paraStart).InsertBefore($fullText)

# position of the bold "Reducer" (the 2nd chunk we concatenated)
$boldReducerStart = $paraStart + $run1.Length
$boldReducerEnd = $boldReducerStart + $run2.Length
$d.Range($boldReducerStart, $boldReducerEnd).Font.Bold = 1

# position of the bold "Store" (the 5th chunk we concatenated)
$boldStoreStart = $paraStart + $run1.Length + $run2.Length + $run3.Length + $run4.Length
$boldStoreEnd = $boldStoreStart + $run5.Length
$d.Range($boldStoreStart, $boldStoreEnd).Font.Bold = 1
